# This script applies updated Betfair Back/Lay odds to sheet1
# of the Jogos_do_Dia workbook, matching the target diff.
# 237 individual cell values are updated across rows 2-18
# (columns F..AO), reflecting refreshed market odds.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.27
$ws.Range("G2").Value = 1.38
$ws.Range("H2").Value = 7.4
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 5.5
$ws.Range("K2").Value = 8.199999999999999
$ws.Range("L2").Value = 1.21
$ws.Range("N2").Value = 3.45
$ws.Range("O2").Value = 1.12
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1.35
$ws.Range("R2").Value = 1.8
$ws.Range("S2").Value = 1.89
$ws.Range("T2").Value = 1.64
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.08
$ws.Range("W2").Value = 3.6
$ws.Range("AN2").Value = 6.8
# Row 3
$ws.Range("F3").Value = 2.12
$ws.Range("I3").Value = 4.6
$ws.Range("K3").Value = 3.25
$ws.Range("L3").Value = 1.58
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 2.54
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 1.52
$ws.Range("Q3").Value = 2.64
$ws.Range("T3").Value = 2.32
$ws.Range("U3").Value = 1.69
$ws.Range("X3").Value = 8.800000000000001
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 30
$ws.Range("AB3").Value = 6.8
$ws.Range("AD3").Value = 21
$ws.Range("AE3").Value = 210
$ws.Range("AH3").Value = 28
$ws.Range("AN3").Value = 30
$ws.Range("AO3").Value = 140
# Row 4
$ws.Range("G4").Value = 1.9
$ws.Range("H4").Value = 4.4
$ws.Range("J4").Value = 3.9
$ws.Range("K4").Value = 5.3
$ws.Range("M4").Value = 1.03
$ws.Range("V4").Value = 1.22
$ws.Range("W4").Value = 2.1
# Row 5
$ws.Range("F5").Value = 1.76
$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 5.5
$ws.Range("I5").Value = 6.2
$ws.Range("J5").Value = 3.55
$ws.Range("K5").Value = 3.8
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 2.74
$ws.Range("O5").Value = 1.47
$ws.Range("P5").Value = 1.6
$ws.Range("Q5").Value = 2.34
$ws.Range("S5").Value = 4.8
$ws.Range("T5").Value = 2.16
$ws.Range("U5").Value = 1.7
$ws.Range("V5").Value = 1.19
$ws.Range("X5").Value = 10
$ws.Range("AB5").Value = 6.8
$ws.Range("AE5").Value = 120
$ws.Range("AH5").Value = 29
$ws.Range("AJ5").Value = 19
$ws.Range("AM5").Value = 220
$ws.Range("AN5").Value = 18
# Row 6
$ws.Range("F6").Value = 3.65
$ws.Range("G6").Value = 4.4
$ws.Range("I6").Value = 2.2
$ws.Range("J6").Value = 3.15
$ws.Range("K6").Value = 3.95
$ws.Range("N6").Value = 3.35
$ws.Range("O6").Value = 1.34
$ws.Range("P6").Value = 1.8
$ws.Range("Q6").Value = 1.98
$ws.Range("R6").Value = 1.3
$ws.Range("S6").Value = 3.6
$ws.Range("U6").Value = 1.98
$ws.Range("W6").Value = 1.3
$ws.Range("AB6").Value = 17
$ws.Range("AC6").Value = 9.800000000000001
# Row 7
$ws.Range("I7").Value = 1.34
$ws.Range("P7").Value = 2.8
$ws.Range("Q7").Value = 1.46
$ws.Range("R7").Value = 1.73
$ws.Range("S7").Value = 2.12
$ws.Range("V7").Value = 3.9
$ws.Range("X7").Value = 36
$ws.Range("Z7").Value = 9.6
$ws.Range("AA7").Value = 11.5
$ws.Range("AC7").Value = 16.5
$ws.Range("AE7").Value = 13
$ws.Range("AO7").Value = 4.1
# Row 8
$ws.Range("F8").Value = 1.98
$ws.Range("G8").Value = 2.12
$ws.Range("H8").Value = 4.6
$ws.Range("I8").Value = 5.8
$ws.Range("J8").Value = 2.84
$ws.Range("K8").Value = 3.5
$ws.Range("M8").Value = 1.09
$ws.Range("N8").Value = 2.88
$ws.Range("O8").Value = 1.37
$ws.Range("P8").Value = 1.65
$ws.Range("Q8").Value = 2
$ws.Range("S8").Value = 3.5
$ws.Range("W8").Value = 1.89
# Row 9
$ws.Range("F9").Value = 2.3
$ws.Range("G9").Value = 2.62
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 3.55
$ws.Range("L9").Value = 1.3
$ws.Range("M9").Value = 1.06
$ws.Range("O9").Value = 1.26
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 1.8
$ws.Range("R9").Value = 1.39
$ws.Range("S9").Value = 3
$ws.Range("T9").Value = 1.64
$ws.Range("U9").Value = 2.22
$ws.Range("V9").Value = 1.4
$ws.Range("W9").Value = 1.62
# Row 10
$ws.Range("F10").Value = 2.84
$ws.Range("G10").Value = 3.3
$ws.Range("H10").Value = 2.4
$ws.Range("L10").Value = 1.39
$ws.Range("N10").Value = 3.55
$ws.Range("O10").Value = 1.31
$ws.Range("P10").Value = 1.89
$ws.Range("Q10").Value = 1.94
$ws.Range("W10").Value = 1.45
$ws.Range("AA10").Value = 40
$ws.Range("AI10").Value = 42
$ws.Range("AJ10").Value = 55
$ws.Range("AL10").Value = 48
# Row 11
$ws.Range("J11").Value = 2.94
$ws.Range("K11").Value = 3.4
$ws.Range("N11").Value = 2.74
# Row 12
$ws.Range("F12").Value = 4.2
$ws.Range("G12").Value = 4.5
$ws.Range("H12").Value = 2.08
$ws.Range("K12").Value = 3.5
$ws.Range("L12").Value = 1.48
$ws.Range("N12").Value = 2.84
$ws.Range("P12").Value = 1.63
$ws.Range("Q12").Value = 2.26
$ws.Range("T12").Value = 1.97
$ws.Range("U12").Value = 1.83
$ws.Range("V12").Value = 1.86
$ws.Range("W12").Value = 1.28
$ws.Range("X12").Value = 10.5
$ws.Range("Z12").Value = 13
$ws.Range("AA12").Value = 34
$ws.Range("AB12").Value = 12.5
$ws.Range("AC12").Value = 8
$ws.Range("AG12").Value = 18
$ws.Range("AH12").Value = 24
$ws.Range("AI12").Value = 55
$ws.Range("AJ12").Value = 110
$ws.Range("AL12").Value = 80
$ws.Range("AM12").Value = 180
$ws.Range("AO12").Value = 23
# Row 13
$ws.Range("F13").Value = 1.55
$ws.Range("G13").Value = 1.65
$ws.Range("I13").Value = 9.6
$ws.Range("J13").Value = 3.9
$ws.Range("K13").Value = 4.4
$ws.Range("L13").Value = 1.38
$ws.Range("M13").Value = 1.09
$ws.Range("N13").Value = 3.15
$ws.Range("P13").Value = 1.73
$ws.Range("Q13").Value = 2.1
$ws.Range("R13").Value = 1.26
$ws.Range("T13").Value = 2.16
$ws.Range("U13").Value = 1.69
$ws.Range("V13").Value = 1.13
$ws.Range("W13").Value = 2.52
$ws.Range("AF13").Value = 10
# Row 14
$ws.Range("G14").Value = 2.48
$ws.Range("T14").Value = 1.92
$ws.Range("U14").Value = 1.88
$ws.Range("W14").Value = 1.67
$ws.Range("AA14").Value = 85
$ws.Range("AC14").Value = 7.6
$ws.Range("AD14").Value = 16.5
$ws.Range("AE14").Value = 55
$ws.Range("AI14").Value = 70
$ws.Range("AO14").Value = 70
# Row 15
$ws.Range("F15").Value = 1.77
$ws.Range("Q15").Value = 1.99
$ws.Range("T15").Value = 1.88
$ws.Range("U15").Value = 1.92
$ws.Range("V15").Value = 1.21
# Row 16
$ws.Range("F16").Value = 3.6
$ws.Range("G16").Value = 4.2
$ws.Range("H16").Value = 2.14
$ws.Range("I16").Value = 2.22
$ws.Range("J16").Value = 3.4
$ws.Range("L16").Value = 1.49
$ws.Range("M16").Value = 1.09
$ws.Range("N16").Value = 2.74
$ws.Range("O16").Value = 1.4
$ws.Range("P16").Value = 1.7
$ws.Range("Q16").Value = 2.32
$ws.Range("T16").Value = 1.9
$ws.Range("U16").Value = 1.89
$ws.Range("V16").Value = 1.81
$ws.Range("W16").Value = 1.32
$ws.Range("X16").Value = 11
$ws.Range("Z16").Value = 12.5
$ws.Range("AC16").Value = 7.6
$ws.Range("AD16").Value = 12.5
$ws.Range("AE16").Value = 28
$ws.Range("AF16").Value = 27
$ws.Range("AG16").Value = 17
$ws.Range("AH16").Value = 22
$ws.Range("AI16").Value = 55
$ws.Range("AJ16").Value = 95
$ws.Range("AK16").Value = 65
$ws.Range("AL16").Value = 75
$ws.Range("AM16").Value = 150
$ws.Range("AN16").Value = 80
# Row 17
$ws.Range("H17").Value = 2.18
$ws.Range("I17").Value = 2.28
$ws.Range("J17").Value = 3.5
$ws.Range("K17").Value = 3.75
$ws.Range("N17").Value = 3.45
$ws.Range("T17").Value = 1.76
$ws.Range("U17").Value = 2.06
$ws.Range("V17").Value = 1.78
$ws.Range("Z17").Value = 13.5
$ws.Range("AE17").Value = 25
$ws.Range("AK17").Value = 48
# Row 18
$ws.Range("I18").Value = 1.54
$ws.Range("T18").Value = 2.22
$ws.Range("U18").Value = 1.64
$ws.Range("V18").Value = 2.8
